# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# El rango de periodos en mora (E16:E51) se reordena de forma ascendente
# (antes estaba en orden descendente) y se actualizan los valores de
# "Valor Mora" (F16:F51) y "Salario Basico" (G16:G51) para reflejar la
# base de datos actualizada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodos en mora, ahora en orden ascendente (jul-2017 .. sep-2009... en
# este caso de 1707 hasta 2009), uno por fila desde la 16 hasta la 51.
$periodos = @(
    "1707","1709","1710",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009"
)

$firstRow = 16
$lastRow = 51

# Salario Basico actualizado: mismo valor para todas las filas del rango.
$nuevoSalarioBasico = 781242

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstRow + $i

    $ws.Range("E$row").Value = $periodos[$i]

    if ($row -le 26) {
        $valorMora = 29509
    } elseif ($row -le 50) {
        $valorMora = 31249
    } else {
        $valorMora = 29166
    }

    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $nuevoSalarioBasico
}
